$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.195.28"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "1.863.28"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'243.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.63%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "'42.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").Value = "'0.2852"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "'0.06473"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").Value = "'20.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.67%  "

$ws.Range("D12").Value = "'0.07694"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.37%  "

$ws.Range("D13").Value = "1.878.16"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "'94.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.43%  "

$ws.Range("D15").Value = "'0.6821"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("D16").Value = "'5.074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("D17").Value = "'268.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "30.189.76"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").Value = "'13.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.65%  "

$ws.Range("D20").Value = "'0.000007536"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("D21").Value = "'0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "2.112.47"
$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").Value = "'0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").Value = "'5.184"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("D25").Value = "'6.105"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("D26").Value = "'9.315"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").Value = "'165.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("D28").Value = "'18.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("D29").Value = "'1.890"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").Value = "'1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.87%  "

$ws.Range("D31").Value = "'0.09812"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("D32").Value = "'1.517"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.04%  "

$ws.Range("D33").Value = "'4.234"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "

$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("D35").Value = "'0.04694"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("D36").Value = "'1.110"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.30%  "

$ws.Range("D37").Value = "'0.6862"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("E38").Value = "  -0.60%  "

$ws.Range("D39").Value = "'0.01847"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").Value = "'2.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.06%  "

$ws.Range("D41").Value = "'6.404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").Value = "'70.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("D43").Value = "'0.9995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").Value = "'0.8367"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "

$ws.Range("D45").Value = "'1.886"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.72%  "

$ws.Range("D46").Value = "'101.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "'0.4064"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.70%  "

$ws.Range("D48").Value = "'9.218"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("E49").Value = "  -2.42%  "

$ws.Range("D50").Value = "'919.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").Value = "'34.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "
